$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data to the right
$ws.Range("A:A").Insert()

# New first column: "Laboratorio" header and "Lab1: Procesos" value
$ws.Range("A1").Value = "Laboratorio"
$ws.Range("A2").Value = "Lab1: Procesos"

# Copy style of B1 (old header style) onto A1 so it matches the other headers
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$ws.Range("A1").Value = "Laboratorio"

# Column widths
# (COM's ColumnWidth quantizes to 1/6-character steps in this host; 19.6667
# lands on the stored width closest to the source file's 20.42578125.)
$ws.Columns.Item(1).ColumnWidth = 19.6667
$ws.Range("B:E").Columns.AutoFit() | Out-Null

$ws.Range("C9").Select()
